$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (error code 112, DECODING_FAILURE): update description to mention the
# newly supported encodings and add a "Postup" (remedy) describing the manual fix.
$ws.Range("C14").Value = "Pre dopytovaný súbor nebolo možné nájsť enkódovanie. Testované sú formáty UTF-16, UTF-8 a windows-1250"
$ws.Range("D14").Value = "Manuálne otvoriť súbor a zmenit enkódovanie na jeden z podporovaných formátov"

# Row 15 (error code 113, UNSUPPORTED_LOG): the remedy now also mentions the .txt extension.
$ws.Range("D15").Value = "Skontrolujte, že názov súboru obsahuje ""KAM"" alebo ""PAP"" a príponu .log resp .txt"

# Column D grew wider to fit the new, longer text.
$ws.Columns.Item(4).ColumnWidth = 75.4

# Move the active selection as recorded by the author when saving.
$ws.Range("C18").Select()
